$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D: "expval"
$ws.Range("D1").Value = "expval"

# Expected-value data for rows 2..102 (column D)
$values = @(
    121.989910896765,
    129.697793924965,
    120.230719344356,
    119.665806146018,
    129.235977741638,
    125.382118497519,
    120.459292128776,
    123.156114267543,
    121.439957713431,
    125.645147056752,
    124.305054076071,
    120.804903848194,
    123.107727170317,
    120.393953132973,
    123.623554782845,
    121.177221026604,
    114.364051053301,
    117.96345785723,
    120.912302756424,
    114.878884940992,
    117.199262416777,
    120.508083398791,
    121.162716989449,
    120.081189000763,
    121.66792732966,
    115.123180283228,
    121.82054305199,
    121.855256265948,
    121.801264282139,
    121.250966284914,
    120.439026342324,
    121.823233207769,
    121.625487586061,
    116.239442104632,
    110.993718510803,
    120.784281215055,
    120.398164905531,
    121.395095285333,
    136.561881387139,
    126.256827797916,
    118.420485667135,
    121.29156295359,
    117.656972745387,
    108.806196223492,
    118.388170912385,
    113.524140591099,
    117.410954622583,
    114.548452373091,
    118.342393709361,
    113.350002529279,
    113.55141763046,
    115.503839879336,
    121.333868445679,
    115.054630258146,
    120.472067984248,
    121.260528503978,
    120.614418054437,
    117.172571155337,
    114.635352739791,
    119.815745223765,
    117.569883216862,
    120.240051719881,
    120.121775493153,
    117.190221183311,
    118.600633960811,
    114.491392538292,
    123.727479773656,
    111.50802726058,
    121.447300012887,
    121.404109754235,
    120.984608342506,
    118.332516098165,
    120.781330384156,
    120.938110164478,
    121.474216090883,
    113.343905012089,
    115.359720949917,
    123.426274049671,
    123.304715164726,
    120.33641442544,
    123.815473958624,
    119.823463989325,
    115.266340194479,
    117.030703376507,
    116.700977382144,
    121.049724109707,
    113.327545443232,
    120.399448200669,
    112.557325807964,
    117.190911231469,
    129.094540919188,
    118.625249125071,
    121.233484609976,
    121.842307787216,
    121.520163542405,
    121.635998526793,
    120.680103700345,
    121.653164044284,
    115.127671719834,
    114.559256285896,
    121.138096580978
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $values[$i]
}

# Column D width (matches the width introduced alongside the new column)
$ws.Columns.Item(4).ColumnWidth = 10.5454545454545

# Selection / view state: extend the selected range to include column D
# and scroll so row 19 is the top visible row (matches the recorded view state).
$ws.Range("A2:D102").Select()
$excel.ActiveWindow.ScrollRow = 19
